$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1").Value = "XlsxTemplate 超链接测试"
}
